$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 9.395728999999999
$ws.Range("I2").Value = 0.780864851881971
$ws.Range("J2").Value = 0.7808648518819711
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.315861666666667
$ws.Range("N2").Value = 3.947585
$ws.Range("O2").Value = 0.2754050739440597
$ws.Range("P2").Value = 0.2754050739440597
$ws.Range("Q2").Value = 4.121159873829445
$ws.Range("R2").Value = 37.090438864465
$ws.Range("S2").Value = 0.2150541422728715
$ws.Range("T2").Value = 0.2150541422728715

# Row 3
$ws.Range("H3").Value = 9.395728999999999
$ws.Range("I3").Value = 0.780864851881971
$ws.Range("J3").Value = 0.7808648518819711
$ws.Range("O3").Value = 0.3040809095127364
$ws.Range("P3").Value = 0.3040809095127364
$ws.Range("Q3").Value = 4.550264905199222
$ws.Range("R3").Value = 40.95238414679299
$ws.Range("S3").Value = 0.2374460943667979
$ws.Range("T3").Value = 0.237446094366798

# Row 4
$ws.Range("H4").Value = 9.395728999999999
$ws.Range("I4").Value = 0.780864851881971
$ws.Range("J4").Value = 0.7808648518819711
$ws.Range("M4").Value = 2.009179666666667
$ws.Range("N4").Value = 6.027539
$ws.Range("O4").Value = 0.4205140165432039
$ws.Range("P4").Value = 0.4205140165432039
$ws.Range("Q4").Value = 6.292569220103445
$ws.Range("R4").Value = 56.633122980931
$ws.Range("S4").Value = 0.3283646152423016
$ws.Range("T4").Value = 0.3283646152423016

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.878912
$ws.Range("H5").Value = 2.636736
$ws.Range("I5").Value = 0.219135148118029
$ws.Range("J5").Value = 0.219135148118029
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.315861666666667
$ws.Range("N5").Value = 3.947585
$ws.Range("O5").Value = 0.2754050739440597
$ws.Range("P5").Value = 0.2754050739440597
$ws.Range("Q5").Value = 1.156526609173333
$ws.Range("R5").Value = 10.40873948256
$ws.Range("S5").Value = 0.06035093167118825
$ws.Range("T5").Value = 0.06035093167118827

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.878912
$ws.Range("H6").Value = 2.636736
$ws.Range("I6").Value = 0.219135148118029
$ws.Range("J6").Value = 0.219135148118029
$ws.Range("O6").Value = 0.3040809095127364
$ws.Range("P6").Value = 0.3040809095127364
$ws.Range("Q6").Value = 1.276946928234667
$ws.Range("R6").Value = 11.492522354112
$ws.Range("S6").Value = 0.06663481514593847
$ws.Range("T6").Value = 0.06663481514593848

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.878912
$ws.Range("H7").Value = 2.636736
$ws.Range("I7").Value = 0.219135148118029
$ws.Range("J7").Value = 0.219135148118029
$ws.Range("M7").Value = 2.009179666666667
$ws.Range("N7").Value = 6.027539
$ws.Range("O7").Value = 0.4205140165432039
$ws.Range("P7").Value = 0.4205140165432039
$ws.Range("Q7").Value = 1.765892119189334
$ws.Range("R7").Value = 15.893029072704
$ws.Range("S7").Value = 0.09214940130090229
$ws.Range("T7").Value = 0.09214940130090229
